$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Handoff xliff generation for 1c06b937/71550e8c/b982bcc0/cc72c133 completed:
# bump the "Latest HO Xliff Generate Date" for those rows, and flip their
# Priority from "low" to "ht" now that they are ready for handoff.

$overview.Range("G4:G7").Value = "2016-08-26 06:29:58"

$zhcn.Range("E4:E7").Value = "ht"
$zhcn.Range("H4:H7").Value = "2016-08-26 06:29:53"

$dede.Range("E4:E7").Value = "ht"
$dede.Range("H4:H7").Value = "2016-08-26 06:29:58"
